# ddepewRubric.xlsx edit: mark model-loader/texture related criteria as met.
# For several rows, set column E (Milestone level, roman numeral) and/or
# column F (the "X" check mark) so the rubric reflects the textures work
# that was completed. Downstream formulas (G, H, I, J, K, L columns) are
# driven by formulas and recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5: mark as met (column E already "III"; check the box in F)
$ws.Range("F5").Value = "X"

# Row 6: mark as met (column E already "III"; check the box in F)
$ws.Range("F6").Value = "X"

# Row 9: mark as met (column E already "III"; check the box in F)
$ws.Range("F9").Value = "X"

# Row 10: mark as met (column E already "III"; check the box in F)
$ws.Range("F10").Value = "X"

# Row 39: set milestone level to III and check the box
$ws.Range("E39").Value = "III"
$ws.Range("F39").Value = "X"

# Row 56: mark as met (column E already "III"; check the box in F)
$ws.Range("F56").Value = "X"

# Update the last-saved view position/selection to reflect where the
# author was working when they made this edit.
$ws.Activate()
$ws.Range("F57").Select()
